$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing data (task now "Pronto", doer "Camila", task text changed)
$ws.Range("A2").Value = "27/03/2024 - 03/04/2024"
$ws.Range("B2").Value = "Criação do Diagrama de Relacionamento"
$ws.Range("C2").Value = "Camila"
$ws.Range("D2").Value = "Pronto"

# B2 gets its own style: centered horizontally (no longer vertically centered)
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4107

# Row 3: new sprint entry
$ws.Range("A3").Value = "27/03/2024 - 03/04/2024"
$ws.Range("B3").Value = "Criação da 1ª versão do Diagrama de UML"
$ws.Range("C3").Value = "Guilherme"
$ws.Range("D3").Value = "Pronto"

# Row 4: new sprint entry
$ws.Range("A4").Value = "27/03/2024 - 03/04/2024"
$ws.Range("B4").Value = "Crud da tela de Home"
$ws.Range("C4").Value = "Bruno e Camila"
$ws.Range("D4").Value = "Pronto"

# New rows follow the same centered style as the rest of the table
$ws.Range("A3:D4").HorizontalAlignment = -4108
$ws.Range("A3:D4").VerticalAlignment = -4108

# Expand the table to include the new rows
$lo = $ws.ListObjects("Tabela2")
$lo.Resize($ws.Range("A1:D4"))

# Move the active selection, as left by the author after editing
[void]$ws.Range("B12").Select()
